$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: add new row 5 ---
$logs.Range("A5").Value = "Tvvccxx hghhgg rtrtrt. Kan dit?"
$logs.Range("B5").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C5").Value = "Hallo, `nTvvccxx hghhgg rtrtrt blabla. Is dit mogelijk op basis van de 7e afmeting van jullie multitemperatuur?`nDank en groet `nJan`nSent using {0}"
$logs.Range("D5").Value = "Productinformatie"
$logs.Range("E5").Value = "Beste Jan,`nBedankt voor je e-mail. Helaas begrijp ik niet helemaal wat je bedoelt met `"Tvvccxx hghhgg rtrtrt blabla`" en `"7e afmeting van jullie multitemperatuur`". Zou je meer details kunnen geven of specifieker kunnen zijn over je vraag?`nAls je nog steeds hulp nodig hebt, laat het me dan alsjeblieft weten. `nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F5").Value = "2025-06-26 18:54:36"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"

# Remove auto-applied custom row height so the row matches the standard height
$logs.Rows.Item(5).AutoFit()

# --- Extend conditional formatting ranges to cover the new row ---
$dFcs = $logs.Range("D2:D4").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D5"))
}

$gFcs = $logs.Range("G2:G4").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G5"))
}

$hFcs = $logs.Range("H2:H4").FormatConditions
for ($i = 1; $i -le $hFcs.Count; $i++) {
    $hFcs.Item($i).ModifyAppliesToRange($logs.Range("H2:H5"))
}

# --- Dashboard sheet: add new row 4 ---
$dashboard.Range("A4").Value = "Productinformatie"
$dashboard.Range("B4").Value = 1

# --- Update chart series references to include the new Dashboard row ---
$chart = $dashboard.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$4"
$series.Values = "='Dashboard'!`$B`$2:`$B`$4"
